$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F27").Value = 3.5
$ws.Range("F28").Value = 5.5
$ws.Range("F29").Value = 6
$ws.Range("F30").Value = 5
$ws.Range("F31").Value = 4
$ws.Range("F32").Value = 1

$ws.Range("H31").Value = "Done"
$ws.Range("H32").Value = "In-Dev"
